# Apply portfolio content updates:
#  - Slide 1: update the "career start" subtitle date
#  - Slide 3: update title/body copy for the "그로스폴리오" project slide

$p = $ppt.ActivePresentation

# --- Slide 1: Subtitle date range ---
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2).TextFrame.TextRange
$subtitle.Paragraphs(1, 1).Runs(1, 1).Text = "2023.04 - 현재"

# --- Slide 3: Title shape with project summary bullets ---
$s3 = $p.Slides.Item(3)
$tr = $s3.Shapes.Item(1).TextFrame.TextRange

$tr.Paragraphs(1, 1).Runs(1, 1).Text = "그로스폴리오 론칭 캠페인"
# Paragraph 2 "프로젝트 상세" unchanged
$tr.Paragraphs(3, 1).Runs(1, 1).Text = "• 신규 서비스 그로스폴리오의 브랜드 아이덴티티 및 마케팅 전략 수립"
$tr.Paragraphs(4, 1).Runs(1, 1).Text = "• 타겟 고객층 분석을 통한 맞춤형 프로모션 캠페인 기획 및 실행"
$tr.Paragraphs(5, 1).Runs(1, 1).Text = "• SNS 채널별 차별화된 콘텐츠 제작 및 운영 전략 수립"
$tr.Paragraphs(6, 1).Runs(1, 1).Text = "• 유저 리텐션 향상을 위한 리워드 프로그램 설계 및 구현"
$tr.Paragraphs(7, 1).Runs(1, 1).Text = "• 서비스 인지도 제고를 위한 디지털 광고 캠페인 운영"
# Paragraph 8 "주요 성과" unchanged
$tr.Paragraphs(9, 1).Runs(1, 1).Text = "• 론칭 3개월 만에 신규 가입자 50,000명 확보"
$tr.Paragraphs(10, 1).Runs(1, 1).Text = "• 캠페인 기간 동안 앱 다운로드 전환율 35% 달성"
$tr.Paragraphs(11, 1).Runs(1, 1).Text = "• 마케팅 활동을 통한 서비스 인지도 22% 상승"
